$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at AG:AH (old AG/AH "Status"/"On/Off" shift right to AI/AJ)
$ws.Columns("AG:AH").Insert()

# Rename the Regional Manager headers to the shorter RM labels
$ws.Range("AE1").Value = "RM Name"
$ws.Range("AF1").Value = "RM Contact No"

# New ASM header + placeholder columns
$ws.Range("AG1").Value = "ASM Name"
$ws.Range("AH1").Value = "ASM Contact No"
$ws.Range("AG2").Value = "{vendor:sf_asm_name}"
$ws.Range("AH2").Value = "{vendor:sf_asm_phone}"

# Match column widths: AG/AH should match AF's width (closest achievable to 33.88671875)
$ws.Columns("AG:AH").ColumnWidth = 33

# Restore the previous selection roughly where the diff left it
$ws.Range("AH8").Select()
